$wb = $excel.ActiveWorkbook

# --- Update status values from "Ready for handoff" to "In Translation" ---

# "zh-cn" worksheet: Status column is C, data rows 2-4
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2:C4").Value = "In Translation"
$wsZh.Columns.Item(3).ColumnWidth = 12.5

# "de-de" worksheet: Status column is C, data rows 2-4
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2:C4").Value = "In Translation"
$wsDe.Columns.Item(3).ColumnWidth = 12.5

# "Overview" worksheet: zh-cn / de-de status columns are E and F, data rows 2-4
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = "In Translation"
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
